# Update marksheet correct/total marks
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("quiz")

# "Marking" row -> Right column (B11): 3 -> 5
$ws.Range("B11").Value = 5

# "Total" row -> Right column (B12): 57 -> 95
$ws.Range("B12").Value = 95

# "Total" row -> Max column (E12): "55/84" -> "95/140"
$ws.Range("E12").Value = "95/140"
